$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 14, pushing existing rows 14-20 down to 15-21
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 45203
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 300000000
$ws.Cells.Item(14, 7).Value = "Espárragos"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 1400
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1450
$ws.Cells.Item(14, 14).Value = "$/kilo"
$ws.Cells.Item(14, 15).Value = "Provincia de Linares"
$ws.Cells.Item(14, 16).Value = 1450
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"
